$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Absent" column (H) to reflect the consolidated report:
# Rows where attendance count was 0 now get an explicit Absent = 1
$ws.Range("H3").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H12").Value = 1

# Rows where attendance was recorded (Real) now get explicit Absent = 0
# instead of being left blank
$ws.Range("H5").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H14").Value = 0
